$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "279.27"
Set-TextValue "E2" "0.54%"
Set-TextValue "G2" "21"
Set-TextValue "D3" "27.47"
Set-TextValue "E3" "1.04%"
Set-TextValue "G3" "21"
Set-TextValue "E4" "-0.73%"
Set-TextValue "G4" "21"
Set-TextValue "D5" "0.06375"
Set-TextValue "E5" "0.20%"
Set-TextValue "G5" "21"
Set-TextValue "D6" "7.019"
Set-TextValue "E6" "0.69%"
Set-TextValue "G6" "21"
Set-TextValue "D7" "1.290"
Set-TextValue "E7" "3.43%"
Set-TextValue "G7" "21"
Set-TextValue "D8" "0.8933"
Set-TextValue "E8" "1.28%"
Set-TextValue "G8" "21"
Set-TextValue "D9" "0.1535"
Set-TextValue "E9" "0.72%"
Set-TextValue "G9" "21"
Set-TextValue "D10" "0.06090"
Set-TextValue "E10" "19.33%"
Set-TextValue "G10" "21"
Set-TextValue "D11" "0.07508"
Set-TextValue "E11" "-0.35%"
Set-TextValue "G11" "21"
Set-TextValue "D12" "0.02937"
Set-TextValue "E12" "-1.20%"
Set-TextValue "G12" "21"
Set-TextValue "D13" "0.08995"
Set-TextValue "E13" "-0.16%"
Set-TextValue "G13" "21"
Set-TextValue "D14" "0.001561"
Set-TextValue "E14" "-0.85%"
Set-TextValue "G14" "21"
Set-TextValue "D15" "0.0006396"
Set-TextValue "E15" "-0.12%"
Set-TextValue "G15" "21"
Set-TextValue "D16" "0.006003"
Set-TextValue "E16" "1.20%"
Set-TextValue "G16" "21"
Set-TextValue "D17" "3.483"
Set-TextValue "E17" "0.61%"
Set-TextValue "G17" "21"
Set-TextValue "D18" "3.325"
Set-TextValue "E18" "0.37%"
Set-TextValue "G18" "21"
Set-TextValue "D19" "2.229"
Set-TextValue "E19" "-1.90%"
Set-TextValue "G19" "21"
Set-TextValue "G20" "21"
Set-TextValue "E21" "1.07%"
Set-TextValue "G21" "21"
Set-TextValue "D22" "3.908"
Set-TextValue "E22" "-0.03%"
Set-TextValue "G22" "21"
Set-TextValue "B23" "ZBToken"
Set-TextValue "C23" "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
Set-TextValue "D23" "0.1503"
Set-TextValue "E23" "8.90%"
Set-TextValue "G23" "21"
Set-TextValue "B24" "CoinExToken"
Set-TextValue "C24" "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue "D24" "0.04400"
Set-TextValue "E24" "-0.62%"
Set-TextValue "G24" "21"
Set-TextValue "E25" "0.19%"
Set-TextValue "G25" "21"
Set-TextValue "D26" "0.004282"
Set-TextValue "E26" "10.62%"
Set-TextValue "G26" "21"
Set-TextValue "G27" "21"
Set-TextValue "E28" "-1.71%"
Set-TextValue "G28" "21"
Set-TextValue "D29" "0.0001653"
Set-TextValue "E29" "-14.64%"
Set-TextValue "G29" "21"
Set-TextValue "G30" "21"
Set-TextValue "G31" "21"
Set-TextValue "G32" "21"
Set-TextValue "G33" "21"
Set-TextValue "G34" "21"
Set-TextValue "G35" "21"
Set-TextValue "G36" "21"
Set-TextValue "G37" "21"
Set-TextValue "G38" "21"
Set-TextValue "G39" "21"
Set-TextValue "D40" "0.04073"
Set-TextValue "E40" "-1.69%"
Set-TextValue "G40" "21"
Set-TextValue "B41" "BKEXToken"
Set-TextValue "C41" "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D41" "0.1400"
Set-TextValue "E41" "18.67%"
Set-TextValue "G41" "21"
Set-TextValue "B42" "KickToken"
Set-TextValue "C42" "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D42" "0.006582"
Set-TextValue "E42" "-3.96%"
Set-TextValue "G42" "21"
Set-TextValue "D43" "0.002078"
Set-TextValue "E43" "2.92%"
Set-TextValue "G43" "21"
Set-TextValue "E44" "-1.95%"
Set-TextValue "G44" "21"
Set-TextValue "D45" "0.00005545"
Set-TextValue "E45" "7.02%"
Set-TextValue "G45" "21"
Set-TextValue "E46" "5.01%"
Set-TextValue "G46" "21"
Set-TextValue "D47" "0.01848"
Set-TextValue "E47" "-8.70%"
Set-TextValue "G47" "21"
Set-TextValue "G48" "21"
Set-TextValue "G49" "21"
Set-TextValue "G50" "21"
Set-TextValue "G51" "21"
